# Fano factor 3d plotting plotly
# Rotate the hypercolumn/layerA/layerB/minicolumn columns (B,C,D,E) one
# position to the right - i.e. new B = old E, new C = old B, new D = old C,
# new E = old D - matching the reordering of columns used to feed the new
# 3d fano-factor plotly surface, and then re-sort the data rows by the new
# (C, B, D, E) key so the table stays grouped for plotting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 37

# --- Read the existing table into memory first ---------------------------
$header = @{}
for ($c = 2; $c -le 5; $c++) {
    $header[$c] = $ws.Cells.Item(1, $c).Value()
}

$data = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $row = [pscustomobject]@{
        A = $ws.Cells.Item($r, 1).Value()
        B = $ws.Cells.Item($r, 2).Value()
        C = $ws.Cells.Item($r, 3).Value()
        D = $ws.Cells.Item($r, 4).Value()
        E = $ws.Cells.Item($r, 5).Value()
    }
    $data += ,$row
}

# --- Rotate header labels: B<-E, C<-B, D<-C, E<-D -------------------------
$newHeader = @{
    2 = $header[5]
    3 = $header[2]
    4 = $header[3]
    5 = $header[4]
}
for ($c = 2; $c -le 5; $c++) {
    $ws.Cells.Item(1, $c).Value = $newHeader[$c]
}

# --- Rotate each data row's columns the same way, building a composite -----
# --- sort key (C, B, D, E) so we can re-sort with a single -Property -------
$rotated = @()
foreach ($row in $data) {
    $newB = $row.E
    $newC = $row.B
    $newD = $row.C
    $newE = $row.D
    $sortKey = "{0:D2}_{1:D2}_{2:D2}_{3:D2}" -f [int]$newC, [int]$newB, [int]$newD, [int]$newE
    $newRow = [pscustomobject]@{
        A   = $row.A
        B   = $newB
        C   = $newC
        D   = $newD
        E   = $newE
        Key = $sortKey
    }
    $rotated += ,$newRow
}

$sorted = $rotated | Sort-Object -Property Key

# --- Write the sorted rows back out ----------------------------------------
$r = 2
foreach ($row in $sorted) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $r++
}
